# Update "Förändrad" (column C) date for all existing data rows (2-44)
# from 45192 to 45202, add a new data row (45) and mark row 44 with an
# explicit row height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the "Förändrad" value in column C for rows 2..44 from 45192 to 45202.
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# 2) Give row 44 an explicit (default) row height, matching the diff.
$ws.Rows.Item(44).RowHeight = 15

# 3) Append the new row 45 with the new cleared-case entry.
$ws.Cells.Item(45, 1).Value = "A 46500-2023"
$ws.Cells.Item(45, 2).Value = 45197
$ws.Cells.Item(45, 3).Value = 45202
$ws.Cells.Item(45, 4).Value = "BLEKINGE LÄN"
$ws.Cells.Item(45, 5).Value = "SÖLVESBORG"
$ws.Cells.Item(45, 7).Value = 0.5
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = 0
$ws.Cells.Item(45, 14).Value = 0
$ws.Cells.Item(45, 15).Value = 0
$ws.Cells.Item(45, 16).Value = 0
$ws.Cells.Item(45, 17).Value = 0

# Match styling used by the other rows: date format for B/C, wrap-text for R.
$ws.Range("B45:C45").NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(45, 18).WrapText = $true
